$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 5.192
$ws.Range("A8").Value = -22.322
$ws.Range("A10").Value = -21.518
$ws.Range("A12").Value = -21.606
$ws.Range("B14").Value = 5.765
$ws.Range("B15").Value = 4.866999999999999
$ws.Range("A18").Value = -21.313
$ws.Range("B18").Value = 7.007
$ws.Range("B20").Value = 6.305
$ws.Range("A25").Value = -21.813
$ws.Range("B29").Value = 5.699
$ws.Range("B30").Value = 5.495000000000001
$ws.Range("B31").Value = 5.751
$ws.Range("B35").Value = 8.059000000000001
$ws.Range("A37").Value = -20.36
$ws.Range("B40").Value = 8.597999999999999
$ws.Range("B44").Value = 5.277
$ws.Range("B50").Value = 4.715000000000001
$ws.Range("B54").Value = 4.995
$ws.Range("A55").Value = -21.797
$ws.Range("A68").Value = -21.507
$ws.Range("B68").Value = 5.881
$ws.Range("B76").Value = 6.343000000000001
$ws.Range("A77").Value = -21.03599999999999
$ws.Range("A78").Value = -20.308
$ws.Range("A79").Value = -21.723
$ws.Range("A80").Value = -20.616
$ws.Range("A81").Value = -21.804
$ws.Range("A82").Value = -22.077
$ws.Range("A84").Value = -21.786
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858
$ws.Range("B92").Value = 6.000999999999999
$ws.Range("B96").Value = 6.427000000000001
$ws.Range("B98").Value = 5.646
$ws.Range("A101").Value = -21.557
$ws.Range("B101").Value = 6.026
$ws.Range("A102").Value = -21.254
$ws.Range("B102").Value = 6.515000000000001
